# Master Data refresh (16th May) - adds 3 new user_detail_h rows (110033-110035)
# for Nikola Tesla, Graham Bell and Albert Miles, mirroring the existing
# data pattern used for rows 2-33 of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up row 33 (its is_active cell was missing the left-aligned style
#     that every other data row in column I already carries) ---
$ws.Range("I33").HorizontalAlignment = -4131   # xlHAlignLeft

# --- New rows appended to the bottom of the table ---
$newRows = @(
    @{ Row = 34; Id = 110033; Uin = 9317596771; Name = "Nikola Tesla"; Email = "nikola.tesla@xyz.com"; Mobile = 818876434 },
    @{ Row = 35; Id = 110034; Uin = 9317596772; Name = "Graham Bell";  Email = "graham.bell@xyz.com";  Mobile = 818876435 },
    @{ Row = 36; Id = 110035; Uin = 9317596773; Name = "Albert Miles"; Email = "albert.miles@xyz.com"; Mobile = 818876436 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Id           # A - id
    $ws.Cells.Item($row, 2).Value = $r.Uin           # B - uin
    $ws.Cells.Item($row, 3).Value = $r.Name          # C - name
    $ws.Cells.Item($row, 4).Value = $r.Email         # D - email
    $ws.Cells.Item($row, 4).HorizontalAlignment = -4131  # xlHAlignLeft (matches rows 2-32)
    $ws.Cells.Item($row, 5).Value = $r.Mobile        # E - mobile
    $ws.Cells.Item($row, 6).Value = "ACT"            # F - status_code
    $ws.Cells.Item($row, 7).Value = "eng"            # G - lang_code
    $ws.Cells.Item($row, 8).Value = "PWD"            # H - last_login_method
    $ws.Cells.Item($row, 9).Value = $true            # I - is_active
    $ws.Cells.Item($row, 9).HorizontalAlignment = -4131  # xlHAlignLeft (matches rows 2-33)
    $ws.Cells.Item($row, 10).Value = "superadmin"    # J - cr_by
    $ws.Cells.Item($row, 11).Value = "now()"         # K - cr_dtimes
    $ws.Cells.Item($row, 12).Value = "now()"         # L - eff_dtimes
}

# --- Restore the selection to the top of the "rest of the sheet" block,
#     matching the refreshed workbook's cursor position ---
$ws.Range("M1:XFD1048576").Select()
